# TC23_Verify_UserRegistration.xlsx
# "Changes done for Javascriptexecutor elements"
#
# Sheet1 ("TC23_Verify_UserRegistration"): the registration flow gained an
# explicit "click register button" step, switched the text-entry clicks to
# CLICK_PRE_ENTERTEXT, and appended a block of CLICK_JS (JavaScript executor)
# steps for the existing-account radio button / account type / T&C checkbox /
# submit button, followed by a success-message verification and a
# my-account/logout sequence.
#
# Sheet2 ("Testdata"): the ExistingaccNObutton data row now points at
# "Profile" rather than "auto", and four new Elementype rows were added for
# the JS-executor steps.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item("Testdata")

# ---------------------------------------------------------------------------
# Sheet1 data (rows 3-29). Column A is always blank in this block.
# ---------------------------------------------------------------------------

# Make sure every cell we are about to touch has the thin all-round border
# that is used throughout this table.
$ws1.Range("A3:E29").Borders.LineStyle = 1

$rows1 = @(
  @{ R=3;  B="WAIT" },
  @{ R=4;  B="CLICK";               C="LoginOption";              D="CSS" },
  @{ R=5;  B="CLICK";               C="RegisterButton";           D="CSS" },
  @{ R=6;  B="CLICK_PRE_ENTERTEXT"; C="RegistrationFName";        D="CSS" },
  @{ R=7;  B="ENTERTEXT";           C="RegistrationFName";        D="CSS"; E="RegistrationFName" },
  @{ R=8;  B="CLICK_PRE_ENTERTEXT"; C="RegistrationLname";        D="CSS" },
  @{ R=9;  B="ENTERTEXT";           C="RegistrationLname";        D="CSS"; E="RegistrationLname" },
  @{ R=10; B="CLICK_PRE_ENTERTEXT"; C="RegistrationEmail";        D="CSS" },
  @{ R=11; B="ENTER_RANDOM_VALUE";  C="RegistrationEmail";        D="CSS"; E="RegistrationEmail" },
  @{ R=12; B="CLICK_PRE_ENTERTEXT"; C="RegistrationPass";         D="CSS" },
  @{ R=13; B="ENTERTEXT";           C="RegistrationPass";         D="CSS"; E="RegistrationPass" },
  @{ R=14; B="SCROLL_DOWN" },
  @{ R=15; B="CLICK_PRE_ENTERTEXT"; C="RegistrationConfirmPass";  D="CSS" },
  @{ R=16; B="ENTERTEXT";           C="RegistrationConfirmPass";  D="CSS"; E="RegistrationConfirmPass" },
  @{ R=17; B="CLICK_PRE_ENTERTEXT"; C="RegistrationZip";          D="CSS" },
  @{ R=18; B="ENTERTEXT";           C="RegistrationZip";          D="CSS"; E="RegistrationZip" },
  @{ R=19; B="CLICK_JS";            C="ExistingaccNObutton";      D="ID";  E="Elementype1" },
  @{ R=20; B="WAIT" },
  @{ R=21; B="CLICK_JS";            C="Acctype";                  D="ID";  E="Elementype2" },
  @{ R=22; B="WAIT" },
  @{ R=23; B="CLICK_JS";            C="T&CCHeckbox";              D="ID";  E="Elementype3" },
  @{ R=24; B="WAIT" },
  @{ R=25; B="CLICK_JS";            C="RegistrationSubmit";       D="ID";  E="Elementype4" },
  @{ R=26; B="WAIT" },
  @{ R=27; B="VERIFY_TEXT_PRESENT"; C="RegistrationSuccessMSG";   D="xpath"; E="RegistrationSuccessMSG" },
  @{ R=28; B="CLICK";               C="MyaccountSection";         D="CSS" },
  @{ R=29; B="CLICK";               C="Logout";                   D="CSS" }
)

foreach ($row in $rows1) {
  $r = $row.R
  $ws1.Cells.Item($r, 2).Value = $row.B
  if ($row.ContainsKey("C")) { $ws1.Cells.Item($r, 3).Value = $row.C } else { $ws1.Cells.Item($r, 3).Value = "" }
  if ($row.ContainsKey("D")) { $ws1.Cells.Item($r, 4).Value = $row.D } else { $ws1.Cells.Item($r, 4).Value = "" }
  if ($row.ContainsKey("E")) { $ws1.Cells.Item($r, 5).Value = $row.E } else { $ws1.Cells.Item($r, 5).Value = "" }
}

# ---------------------------------------------------------------------------
# Sheet2 data (Testdata). Rows 1-8 are unchanged; row 9's value changes and
# four new rows are appended.
# ---------------------------------------------------------------------------

$ws2.Range("A9:B13").Borders.LineStyle = 1
$ws2.Range("B9").Value = "Profile"

$rows2 = @(
  @{ R=10; A="Elementype1"; B="RadioButton" },
  @{ R=11; A="Elementype2"; B="RadioButton" },
  @{ R=12; A="Elementype3"; B="Checkbox" },
  @{ R=13; A="Elementype4"; B="Button" }
)

foreach ($row in $rows2) {
  $r = $row.R
  $ws2.Cells.Item($r, 1).Value = $row.A
  $ws2.Cells.Item($r, 2).Value = $row.B
}

# ---------------------------------------------------------------------------
# View state: keep sheet1 as the active/selected tab with its selection on
# the newly appended block, and leave sheet2's selection on its new rows.
# ---------------------------------------------------------------------------

$ws2.Range("A9:B13").Select()
$ws1.Activate()
$ws1.Range("A27:XFD28").Select()
